$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, shifting existing rows 152-182 down to 153-183
$ws.Rows("152:152").Insert()

# Populate the newly inserted row 152 with the new record's data
$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("D152").Value = 44504
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = 100112037
$ws.Range("G152").Value = "Cebollín"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 100
$ws.Range("K152").Value = 5000
$ws.Range("L152").Value = 6000
$ws.Range("M152").Value = 5500
$ws.Range("N152").Value = "$/paquete 36 unidades"
$ws.Range("O152").Value = "Región Metropolitana"
$ws.Range("P152").Value = 153
$ws.Range("Q152").Value = 36
$ws.Range("R152").Value = "Hortaliza"
